$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for two new rows right after "Spiral5" (row 3) ---
# The existing "RotRing OmegaMax-90" ... "Michael-SNHex" block (rows 4:29)
# needs to move down to rows 6:31 to make space for the two new entries
# "Holden" and "Rizzie Spiral". Columns C:W already hold a uniform literal
# 1 on every data row, so only column B (label) needs moving; shift
# bottom-up so we never clobber a row before it has been read. Column A
# is just the running 0-based index (row - 2), so it is set directly.
for ($r = 29; $r -ge 4; $r--) {
    $dest = $r + 2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 1).Value = $dest - 2
}

# Rows 30:31 are brand new cells beyond the sheet's old extent (which
# stopped at row 29), so column A there needs the same bold/bordered
# look as every other row's index cell - copy it from an existing row.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)

# --- Fill in the two newly freed rows ---
# Row 4: index 2, label "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"

# Row 5: index 3, label "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"

# Every data row's C:W holds a uniform literal 1 (incl. the two brand-new
# rows above and the two rows now freshly occupied at the bottom, 30:31).
$ws.Range("C4:W31").Value = 1

# --- Rename "Thomas Hex" -> "Matthies Hex" (now sitting in row 11 after the shift) ---
$ws.Range("B11").Value = "Matthies Hex"
